$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the R10 rule (cell E8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Make the sheet active and select/activate cell E8, matching the saved selection state
$ws.Activate()
$ws.Range("E8").Select()
